$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("hdhewgeyuf")

$ws.Cells.Item(1, 10).Value = 32.89208173751831
$ws.Cells.Item(2, 2).Value = 1864
$ws.Cells.Item(2, 4).Value = 1863
$ws.Cells.Item(2, 5).Value = 0
$ws.Cells.Item(2, 6).Value = 1
$ws.Cells.Item(2, 7).Value = 99.94635193133047
$ws.Cells.Item(2, 8).Value = 100
$ws.Cells.Item(2, 9).Value = 0.0005361930294906167
$ws.Cells.Item(2, 10).Value = 43.58099246025085
$ws.Cells.Item(3, 10).Value = 39.80966877937317
$ws.Cells.Item(4, 2).Value = 2601
$ws.Cells.Item(4, 4).Value = 2566
$ws.Cells.Item(4, 6).Value = 5
$ws.Cells.Item(4, 7).Value = 99.80552314274601
$ws.Cells.Item(4, 8).Value = 98.69230769230769
$ws.Cells.Item(4, 9).Value = 0.01516329704510109
$ws.Cells.Item(4, 10).Value = 34.26331782341003
$ws.Cells.Item(5, 2).Value = 2025
$ws.Cells.Item(5, 5).Value = 0
$ws.Cells.Item(5, 8).Value = 100
$ws.Cells.Item(5, 9).Value = 0.0009866798223976321
$ws.Cells.Item(5, 10).Value = 40.78151631355286
$ws.Cells.Item(6, 2).Value = 1731
$ws.Cells.Item(6, 4).Value = 1726
$ws.Cells.Item(6, 5).Value = 4
$ws.Cells.Item(6, 6).Value = 36
$ws.Cells.Item(6, 7).Value = 97.95686719636777
$ws.Cells.Item(6, 8).Value = 99.76878612716763
$ws.Cells.Item(6, 9).Value = 0.02268859897901305
$ws.Cells.Item(6, 10).Value = 35.25825715065002
$ws.Cells.Item(7, 10).Value = 39.45233535766602
$ws.Cells.Item(8, 10).Value = 33.24770474433899
$ws.Cells.Item(9, 10).Value = 37.23942470550537
$ws.Cells.Item(10, 2).Value = 1822
$ws.Cells.Item(10, 5).Value = 28
$ws.Cells.Item(10, 8).Value = 98.46238330587589
$ws.Cells.Item(10, 9).Value = 0.01615598885793872
$ws.Cells.Item(10, 10).Value = 39.76329469680786
$ws.Cells.Item(11, 2).Value = 1800
$ws.Cells.Item(11, 4).Value = 1798
$ws.Cells.Item(11, 5).Value = 1
$ws.Cells.Item(11, 6).Value = 80
$ws.Cells.Item(11, 7).Value = 95.74014909478169
$ws.Cells.Item(11, 8).Value = 99.94441356309061
$ws.Cells.Item(11, 9).Value = 0.04310803618946248
$ws.Cells.Item(11, 10).Value = 29.73467016220093
$ws.Cells.Item(12, 10).Value = 33.30852890014648
$ws.Cells.Item(13, 2).Value = 2389
$ws.Cells.Item(13, 5).Value = 0
$ws.Cells.Item(13, 8).Value = 100
$ws.Cells.Item(13, 9).Value = 0.009535655058043118
$ws.Cells.Item(13, 10).Value = 31.73519277572632
$ws.Cells.Item(14, 2).Value = 1535
$ws.Cells.Item(14, 5).Value = 0
$ws.Cells.Item(14, 8).Value = 100
$ws.Cells.Item(14, 9).Value = 0
$ws.Cells.Item(14, 10).Value = 31.41628885269165
$ws.Cells.Item(15, 2).Value = 2284
$ws.Cells.Item(15, 5).Value = 6
$ws.Cells.Item(15, 8).Value = 99.73718791064388
$ws.Cells.Item(15, 9).Value = 0.002633889376646181
$ws.Cells.Item(15, 10).Value = 32.93444752693176
$ws.Cells.Item(16, 2).Value = 1991
$ws.Cells.Item(16, 5).Value = 4
$ws.Cells.Item(16, 8).Value = 99.79899497487438
$ws.Cells.Item(16, 9).Value = 0.002013085052843483
$ws.Cells.Item(16, 10).Value = 35.55498957633972
$ws.Cells.Item(17, 10).Value = 36.79891228675842
$ws.Cells.Item(18, 10).Value = 35.36628532409668
$ws.Cells.Item(19, 2).Value = 1518
$ws.Cells.Item(19, 5).Value = 0
$ws.Cells.Item(19, 8).Value = 100
$ws.Cells.Item(19, 9).Value = 0
$ws.Cells.Item(19, 10).Value = 35.79813051223755
$ws.Cells.Item(20, 2).Value = 1613
$ws.Cells.Item(20, 5).Value = 0
$ws.Cells.Item(20, 8).Value = 100
$ws.Cells.Item(20, 9).Value = 0.003705991352686844
$ws.Cells.Item(20, 10).Value = 28.27710938453674
$ws.Cells.Item(21, 2).Value = 2602
$ws.Cells.Item(21, 4).Value = 2597
$ws.Cells.Item(21, 5).Value = 4
$ws.Cells.Item(21, 6).Value = 3
$ws.Cells.Item(21, 7).Value = 99.88461538461539
$ws.Cells.Item(21, 8).Value = 99.84621299500192
$ws.Cells.Item(21, 9).Value = 0.002691272587466359
$ws.Cells.Item(21, 10).Value = 33.82521653175354
$ws.Cells.Item(22, 2).Value = 1935
$ws.Cells.Item(22, 4).Value = 1934
$ws.Cells.Item(22, 6).Value = 28
$ws.Cells.Item(22, 7).Value = 98.57288481141693
$ws.Cells.Item(22, 9).Value = 0.01426388181355069
$ws.Cells.Item(22, 10).Value = 39.73013257980347
$ws.Cells.Item(23, 2).Value = 2134
$ws.Cells.Item(23, 4).Value = 2133
$ws.Cells.Item(23, 5).Value = 0
$ws.Cells.Item(23, 6).Value = 2
$ws.Cells.Item(23, 7).Value = 99.90632318501171
$ws.Cells.Item(23, 8).Value = 100
$ws.Cells.Item(23, 10).Value = 30.36996459960938
$ws.Cells.Item(24, 2).Value = 2983
$ws.Cells.Item(24, 4).Value = 2962
$ws.Cells.Item(24, 5).Value = 20
$ws.Cells.Item(24, 6).Value = 17
$ws.Cells.Item(24, 7).Value = 99.42933870426317
$ws.Cells.Item(24, 8).Value = 99.32930918846412
$ws.Cells.Item(24, 9).Value = 0.01241610738255034
$ws.Cells.Item(24, 10).Value = 35.44391465187073
$ws.Cells.Item(25, 2).Value = 2647
$ws.Cells.Item(25, 4).Value = 2646
$ws.Cells.Item(25, 6).Value = 9
$ws.Cells.Item(25, 7).Value = 99.66101694915254
$ws.Cells.Item(25, 9).Value = 0.00338855421686747
$ws.Cells.Item(25, 10).Value = 35.04801058769226
$ws.Cells.Item(26, 2).Value = 1835
$ws.Cells.Item(26, 4).Value = 1830
$ws.Cells.Item(26, 5).Value = 4
$ws.Cells.Item(26, 6).Value = 29
$ws.Cells.Item(26, 7).Value = 98.4400215169446
$ws.Cells.Item(26, 8).Value = 99.78189749182116
$ws.Cells.Item(26, 9).Value = 0.01774193548387097
$ws.Cells.Item(26, 10).Value = 32.77532434463501
$ws.Cells.Item(27, 2).Value = 2945
$ws.Cells.Item(27, 5).Value = 4
$ws.Cells.Item(27, 8).Value = 99.86413043478261
$ws.Cells.Item(27, 9).Value = 0.006091370558375634
$ws.Cells.Item(27, 10).Value = 34.08404612541199
$ws.Cells.Item(28, 2).Value = 3007
$ws.Cells.Item(28, 5).Value = 2
$ws.Cells.Item(28, 8).Value = 99.93346640053227
$ws.Cells.Item(28, 9).Value = 0.0006655574043261231
$ws.Cells.Item(28, 10).Value = 33.1038191318512
$ws.Cells.Item(29, 2).Value = 2650
$ws.Cells.Item(29, 4).Value = 2643
$ws.Cells.Item(29, 5).Value = 6
$ws.Cells.Item(29, 6).Value = 6
$ws.Cells.Item(29, 7).Value = 99.77349943374858
$ws.Cells.Item(29, 8).Value = 99.77349943374858
$ws.Cells.Item(29, 9).Value = 0.004528301886792453
$ws.Cells.Item(29, 10).Value = 35.34947443008423
$ws.Cells.Item(30, 2).Value = 2750
$ws.Cells.Item(30, 5).Value = 2
$ws.Cells.Item(30, 8).Value = 99.9272462713714
$ws.Cells.Item(30, 9).Value = 0.000727802037845706
$ws.Cells.Item(30, 10).Value = 35.17819762229919
$ws.Cells.Item(31, 10).Value = 34.5469434261322
$ws.Cells.Item(32, 2).Value = 2261
$ws.Cells.Item(32, 5).Value = 2
$ws.Cells.Item(32, 8).Value = 99.91150442477876
$ws.Cells.Item(32, 9).Value = 0.002210433244916004
$ws.Cells.Item(32, 10).Value = 36.6521327495575
$ws.Cells.Item(33, 10).Value = 34.40062165260315
$ws.Cells.Item(34, 10).Value = 37.54677772521973
$ws.Cells.Item(35, 10).Value = 43.30392718315125
$ws.Cells.Item(36, 2).Value = 2425
$ws.Cells.Item(36, 4).Value = 2424
$ws.Cells.Item(36, 5).Value = 0
$ws.Cells.Item(36, 6).Value = 2
$ws.Cells.Item(36, 7).Value = 99.91755976916735
$ws.Cells.Item(36, 8).Value = 100
$ws.Cells.Item(36, 9).Value = 0.0008240626287597857
$ws.Cells.Item(36, 10).Value = 34.48687863349915
$ws.Cells.Item(37, 2).Value = 2477
$ws.Cells.Item(37, 4).Value = 2475
$ws.Cells.Item(37, 5).Value = 1
$ws.Cells.Item(37, 6).Value = 7
$ws.Cells.Item(37, 7).Value = 99.71796937953263
$ws.Cells.Item(37, 8).Value = 99.95961227786754
$ws.Cells.Item(37, 9).Value = 0.003221908981071285
$ws.Cells.Item(37, 10).Value = 34.89000248908997
$ws.Cells.Item(38, 10).Value = 29.9749174118042
$ws.Cells.Item(39, 2).Value = 2070
$ws.Cells.Item(39, 4).Value = 2046
$ws.Cells.Item(39, 5).Value = 23
$ws.Cells.Item(39, 6).Value = 6
$ws.Cells.Item(39, 7).Value = 99.70760233918129
$ws.Cells.Item(39, 8).Value = 98.88835186080232
$ws.Cells.Item(39, 9).Value = 0.01412566975158305
$ws.Cells.Item(39, 10).Value = 37.55591177940369
$ws.Cells.Item(40, 10).Value = 39.00035047531128
$ws.Cells.Item(41, 10).Value = 35.43409395217896
$ws.Cells.Item(42, 2).Value = 1780
$ws.Cells.Item(42, 5).Value = 0
$ws.Cells.Item(42, 8).Value = 100
$ws.Cells.Item(42, 9).Value = 0
$ws.Cells.Item(42, 10).Value = 34.44174456596375
$ws.Cells.Item(43, 10).Value = 39.73002004623413
$ws.Cells.Item(44, 10).Value = 34.96149754524231

$ws.Name = "fdcryvy"
